$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dimensions "nivel-estudios" and "sexo" have been re-curated from
# dimensions into measures, so update the descriptor rows accordingly.

# Row 2: semantic identifier changes from iaest-dimension:* to iaest-measure:*
$ws.Range("E2").Value = "iaest-measure:nivel-estudios"
$ws.Range("F2").Value = "iaest-measure:sexo"

# Row 3: classification changes from "dim" to "medida"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "medida"

# Row 4: datatype changes from "skos:Concept" to "xsd:int"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:int"

# Row 5 (mapping file references) is no longer needed and is removed entirely.
$ws.Rows.Item(5).Delete()
